# [ADD] New normalize way
# Swap the "Trafico de drogas" and "Orden publico" columns (E and F),
# including the header row and all data rows, to reflect the new
# normalization ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$eRange = $ws.Range("E1:E51")
$fRange = $ws.Range("F1:F51")

$eVals = $eRange.Value2
$fVals = $fRange.Value2

$eRange.Value2 = $fVals
$fRange.Value2 = $eVals
